$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") rows 2-76 all change from serial date 45175 to 45177
$ws.Range("C2:C76").Value = 45177
